$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Clone File") holds the clone-file names with a ".xls" extension.
# Rename them to the ".xlsx" extension for every data row (2-9).
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 -replace '\.xls$', '.xlsx'
}

# Move the active selection to E10, matching the saved workbook state.
$ws.Range("E10").Select() | Out-Null
